# "atualizado local do projeto, aguardando planilhas a serem consolidadas"
#
# The report sheet held a single placeholder cycle (IN0006 repeated for
# rows 2-4). The consolidated data is now in: a 4-line item cycle
# (IN0006 / IN0012 / IN0018 / IN0025) repeated 3x across rows 2-13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$items = @(
    @{ Id = "IN0006"; Item = "Item 6";  Descr = "Descr 6";  Flag = 1; F = 11; G = 5;  H = 55;  I = 9;  J = 13; K = 150 },
    @{ Id = "IN0012"; Item = "Item 12"; Descr = "Descr 12"; Flag = 1; F = 18; G = 22; H = 396; I = 36; J = 12; K = 50 },
    @{ Id = "IN0018"; Item = "Item 18"; Descr = "Descr 18"; Flag = 1; F = 12; G = 6;  H = 72;  I = 7;  J = 13; K = 50 },
    @{ Id = "IN0025"; Item = "Item 25"; Descr = "Descr 25"; Flag = 0; F = 14; G = 28; H = 392; I = 21; J = 8;  K = 50 }
)

$row = 2
for ($cycle = 0; $cycle -lt 3; $cycle++) {
    for ($idx = 0; $idx -lt $items.Count; $idx++) {
        $data = $items[$idx]
        $ws.Cells.Item($row, 2).Value = $data.Flag
        $ws.Cells.Item($row, 3).Value = $data.Id
        $ws.Cells.Item($row, 4).Value = $data.Item
        $ws.Cells.Item($row, 5).Value = $data.Descr
        $ws.Cells.Item($row, 6).Value = $data.F
        $ws.Cells.Item($row, 7).Value = $data.G
        $ws.Cells.Item($row, 8).Value = $data.H
        $ws.Cells.Item($row, 9).Value = $data.I
        $ws.Cells.Item($row, 10).Value = $data.J
        $ws.Cells.Item($row, 11).Value = $data.K

        $row++
    }
}
